$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2159.742
$ws.Range("I43").Value = 2345.4546
$ws.Range("K43").Value = 2345.4546
$ws.Range("M43").Value = -2276.4546

$ws.Range("H58").Value = 1037.875
$ws.Range("I58").Value = 389.55554
$ws.Range("J58").Value = 1871.4286
$ws.Range("K58").Value = 1168.66662
$ws.Range("L58").Value = 5614.2858
$ws.Range("M58").Value = -1018.66662
$ws.Range("N58").Value = -5914.2858

$ws.Range("H94").Value = 4724.9443
$ws.Range("I94").Value = 4724.9443
$ws.Range("K94").Value = 4724.9443
$ws.Range("M94").Value = -4273.9443

$ws.Range("H108").Value = 26860.334
$ws.Range("J108").Value = 26860.334
$ws.Range("L108").Value = 26860.334
$ws.Range("N108").Value = -34540.334

$ws.Range("H109").Value = 39171
$ws.Range("J109").Value = 39171
$ws.Range("L109").Value = 39171
$ws.Range("N109").Value = -41945

$ws.Range("H120").Value = 49644
$ws.Range("J120").Value = 49644
$ws.Range("L120").Value = 49644
$ws.Range("N120").Value = -59320

$ws.Range("H130").Value = 48830
$ws.Range("J130").Value = 48830
$ws.Range("L130").Value = 48830
$ws.Range("N130").Value = -58870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 58667.11
$ws.Range("I6").Value = 35600.4
$ws.Range("K6").Value = 35600.4
$ws.Range("M6").Value = -35427.4

$ws.Range("H41").Value = 3428
$ws.Range("I41").Value = 3428
$ws.Range("K41").Value = 3428
$ws.Range("M41").Value = -3014

$ws.Range("H61").Value = 3368.8
$ws.Range("I61").Value = 2076.4285
$ws.Range("J61").Value = 4499.625
$ws.Range("K61").Value = 2076.4285
$ws.Range("L61").Value = 4499.625
$ws.Range("M61").Value = -1864.4285
$ws.Range("N61").Value = -4923.625

$ws.Range("H74").Value = 1736.4375
$ws.Range("I74").Value = 698.63635
$ws.Range("J74").Value = 4019.6
$ws.Range("K74").Value = 698.63635
$ws.Range("L74").Value = 4019.6
$ws.Range("M74").Value = 175.36365
$ws.Range("N74").Value = -5767.6

$ws.Range("H77").Value = 1736.4375
$ws.Range("I77").Value = 698.63635
$ws.Range("J77").Value = 4019.6
$ws.Range("K77").Value = 3493.18175
$ws.Range("L77").Value = 20098
$ws.Range("M77").Value = 874.8182500000003
$ws.Range("N77").Value = -28834

$ws.Range("H97").Value = 1309.5714
$ws.Range("J97").Value = 1335.0834
$ws.Range("L97").Value = 1335.0834
$ws.Range("N97").Value = -2327.0834

$ws.Range("H119").Value = 52694
$ws.Range("J119").Value = 52694
$ws.Range("L119").Value = 52694
$ws.Range("N119").Value = -62370

$ws.Range("H133").Value = 33348.418
$ws.Range("J133").Value = 33348.418
$ws.Range("L133").Value = 33348.418
$ws.Range("N133").Value = -38408.418

$ws.Range("H134").Value = 50739.8
$ws.Range("J134").Value = 50739.8
$ws.Range("L134").Value = 50739.8
$ws.Range("N134").Value = -60879.8

$ws.Range("H136").Value = 3368.8
$ws.Range("I136").Value = 2076.4285
$ws.Range("J136").Value = 4499.625
$ws.Range("K136").Value = 6229.2855
$ws.Range("L136").Value = 13498.875
$ws.Range("M136").Value = -3679.2855
$ws.Range("N136").Value = -18598.875

$ws.Range("H138").Value = 47332.25
$ws.Range("J138").Value = 47332.25
$ws.Range("L138").Value = 47332.25
$ws.Range("N138").Value = -57612.25

$ws.Range("H140").Value = 39499.332
$ws.Range("J140").Value = 39499.332
$ws.Range("L140").Value = 39499.332
$ws.Range("N140").Value = -49859.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

$ws.Range("H130").Value = 38274.125
$ws.Range("J130").Value = 38274.125
$ws.Range("L130").Value = 38274.125
$ws.Range("N130").Value = -48314.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 226943.5
$ws.Range("J4").Value = 226943.5
$ws.Range("L4").Value = 226943.5
$ws.Range("N4").Value = -227167.5

$ws.Range("H38").Value = 45000
$ws.Range("J38").Value = 45000
$ws.Range("L38").Value = 45000
$ws.Range("N38").Value = -45754

$ws.Range("H46").Value = 45000
$ws.Range("J46").Value = 45000
$ws.Range("L46").Value = 45000
$ws.Range("N46").Value = -45422

$ws.Range("H138").Value = 45525.668
$ws.Range("J138").Value = 45525.668
$ws.Range("L138").Value = 45525.668
$ws.Range("N138").Value = -55805.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 1750
$ws.Range("J42").Value = 571.4286
$ws.Range("L42").Value = 1714.2858
$ws.Range("N42").Value = -2782.2858

$ws.Range("H68").Value = 6640.6
$ws.Range("J68").Value = 10534.333
$ws.Range("L68").Value = 31602.999
$ws.Range("N68").Value = -33224.999

$ws.Range("H71").Value = 6640.6
$ws.Range("J71").Value = 10534.333
$ws.Range("L71").Value = 94808.997
$ws.Range("N71").Value = -102920.997

$ws.Range("H113").Value = 7901.857
$ws.Range("I113").Value = 17471.834
$ws.Range("J113").Value = 724.375
$ws.Range("K113").Value = 52415.50199999999
$ws.Range("L113").Value = 2173.125
$ws.Range("M113").Value = -50245.50199999999
$ws.Range("N113").Value = -6513.125

$ws.Range("H127").Value = 983.25
$ws.Range("J127").Value = 983.25
$ws.Range("L127").Value = 2949.75
$ws.Range("N127").Value = -12869.75

$ws.Range("H133").Value = 7240.5386
$ws.Range("I133").Value = 7288.3335
$ws.Range("J133").Value = 7199.5713
$ws.Range("K133").Value = 21865.0005
$ws.Range("L133").Value = 21598.7139
$ws.Range("M133").Value = -16805.0005
$ws.Range("N133").Value = -31718.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 22344
$ws.Range("J26").Value = 22344
$ws.Range("L26").Value = 22344
$ws.Range("N26").Value = -22904

$ws.Range("H50").Value = 22344
$ws.Range("J50").Value = 22344
$ws.Range("L50").Value = 22344
$ws.Range("N50").Value = -23340

$ws.Range("H116").Value = 39000
$ws.Range("J116").Value = 39000
$ws.Range("L116").Value = 39000
$ws.Range("N116").Value = -48178

$ws.Range("H119").Value = 27918
$ws.Range("J119").Value = 27918
$ws.Range("L119").Value = 27918
$ws.Range("N119").Value = -37594

$ws.Range("H130").Value = 46312.332
$ws.Range("J130").Value = 46312.332
$ws.Range("L130").Value = 46312.332
$ws.Range("N130").Value = -56352.332

$ws.Range("H135").Value = 30986.375
$ws.Range("J135").Value = 30986.375
$ws.Range("L135").Value = 30986.375
$ws.Range("N135").Value = -41126.375

$ws.Range("H138").Value = 43808.332
$ws.Range("J138").Value = 43808.332
$ws.Range("L138").Value = 43808.332
$ws.Range("N138").Value = -54088.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 24318.273

$ws.Range("H93").Value = 1446.28
$ws.Range("I93").Value = 1071.0834
$ws.Range("J93").Value = 1792.6154
$ws.Range("K93").Value = 1071.0834
$ws.Range("L93").Value = 1792.6154
$ws.Range("M93").Value = 176.9166
$ws.Range("N93").Value = -4288.6154

$ws.Range("H100").Value = 1675.1
$ws.Range("I100").Value = 1593.875
$ws.Range("K100").Value = 1593.875
$ws.Range("M100").Value = -1052.875

$ws.Range("H121").Value = 16905.5
$ws.Range("J121").Value = 16905.5
$ws.Range("L121").Value = 16905.5
$ws.Range("N121").Value = -20399.5

$ws.Range("H127").Value = 44336.5
$ws.Range("J127").Value = 44336.5
$ws.Range("L127").Value = 44336.5
$ws.Range("N127").Value = -54256.5

$ws.Range("H136").Value = 2447.08
$ws.Range("I136").Value = 1688.7368
$ws.Range("J136").Value = 4848.5
$ws.Range("K136").Value = 5066.2104
$ws.Range("L136").Value = 14545.5
$ws.Range("M136").Value = -2516.2104
$ws.Range("N136").Value = -19645.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4167.3335

$ws.Range("H4").Value = 1333.3334
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 1950
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 1950
$ws.Range("M4").Value = 13
$ws.Range("N4").Value = -2176

$ws.Range("H16").Value = 40988.668
$ws.Range("J16").Value = 40988.668
$ws.Range("L16").Value = 40988.668
$ws.Range("N16").Value = -41572.668

$ws.Range("H81").Value = 1457.7778
$ws.Range("I81").Value = 1374.2858
$ws.Range("J81").Value = 1750
$ws.Range("K81").Value = 2748.5716
$ws.Range("L81").Value = 3500
$ws.Range("M81").Value = -1687.5716
$ws.Range("N81").Value = -5622

$ws.Range("H84").Value = 1457.7778
$ws.Range("I84").Value = 1374.2858
$ws.Range("J84").Value = 1750
$ws.Range("K84").Value = 13742.858
$ws.Range("L84").Value = 17500
$ws.Range("M84").Value = -8438.858
$ws.Range("N84").Value = -28108

$ws.Range("H110").Value = 25104
$ws.Range("J110").Value = 25104
$ws.Range("L110").Value = 25104
$ws.Range("N110").Value = -33284

$ws.Range("H136").Value = 21843
$ws.Range("I136").Value = 35222.656
$ws.Range("J136").Value = 2442.5
$ws.Range("K136").Value = 105667.968
$ws.Range("L136").Value = 7327.5
$ws.Range("M136").Value = -103117.968
$ws.Range("N136").Value = -12427.5

$ws.Range("H137").Value = 36480
$ws.Range("J137").Value = 36480
$ws.Range("L137").Value = 36480
$ws.Range("N137").Value = -46680
